$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "index" column (A) data cells down
# into the eight new rows so the new A26:A33 cells pick up the same style
# (bold, centered, thin border) as A2:A25.
$ws.Range("A25").Copy()
$ws.Range("A26:A33").PasteSpecial(-4122)

# New block of rows for the "transcription regulator" GO group, mirroring
# the shape of the previous four 8-row blocks (1-mer .. 8-mer).
$newRows = @(
    @{ Row = 26; A = 0; B = "1-mer"; C = 1; D = 5;  E = 0.108695652173913 },
    @{ Row = 27; A = 1; B = "2-mer"; C = 2; D = 31; E = 0.6739130434782609 },
    @{ Row = 28; A = 2; B = "3-mer"; C = 3; D = 0;  E = 0 },
    @{ Row = 29; A = 3; B = "4-mer"; C = 4; D = 7;  E = 0.1521739130434783 },
    @{ Row = 30; A = 4; B = "5-mer"; C = 5; D = 0;  E = 0 },
    @{ Row = 31; A = 5; B = "6-mer"; C = 6; D = 2;  E = 0.04347826086956522 },
    @{ Row = 32; A = 6; B = "7-mer"; C = 7; D = 0;  E = 0 },
    @{ Row = 33; A = 7; B = "8-mer"; C = 8; D = 1;  E = 0.02173913043478261 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = "transcription regulator"
    $ws.Range("G$row").Value = "Molecular function"
}
